$d = $word.ActiveDocument

# Locate the insertion point: right after "...One-Hot Encoding. "
# (the sentence that previously continued straight into "The original columns...").
$rng = $d.Content
[void]$rng.Find.Execute("Categorical values are converted into dummy variables using One-Hot Encoding. ")
$ip = $rng.Duplicate
$ip.Collapse(0)

# Insert the new sentence about dropping the first dummy value as three
# separate runs (mirroring how the edit was authored), remembering their
# start/end offsets so we can color them afterwards.
$ip.InsertAfter("For each dummy ")
$r1s = $ip.Start
$r1e = $ip.End
$ip.Collapse(0)

$ip.InsertAfter("we")
$r2s = $ip.Start
$r2e = $ip.End
$ip.Collapse(0)

$ip.InsertAfter(" dropped the first value, to prevent multicollinearity")
$r3s = $ip.Start
$r3e = $ip.End
$ip.Collapse(0)

# Close out the new sentence in the normal (non-highlighted) color.
$ip.InsertAfter(". ")
$ip.Collapse(0)

# The remaining original text ("The original columns of ...") now
# immediately follows. Split "The" off into its own run, matching the
# structure of the authored edit, without altering its text/formatting.
$theRng = $d.Range($ip.End, $ip.End + 3)
$theRng.Bold = 1
$theRng.Bold = 0

# Color the three new runs red (FF0000), from last to first so each
# InsertAfter'd run keeps its own identity instead of merging back with
# its already-colored neighbor.
$d.Range($r3s, $r3e).Font.Color = 255
$d.Range($r2s, $r2e).Font.Color = 255
$d.Range($r1s, $r1e).Font.Color = 255
